$wb = $excel.ActiveWorkbook

# =====================================================================
# Discount sheet ("Discount" / sheet7.xml) — value corrections on rows
# 2-7 plus a brand-new row 8 ("Auto 7").
# =====================================================================
$ws7 = $wb.Worksheets.Item("Discount")
$ws7.Activate()

# Row 2
$ws7.Range("C2").Value = 10
$ws7.Range("E2").Value = "yes"

# Row 3
$ws7.Range("B3").Value = 25
$ws7.Range("C3").Value = 25
$ws7.Range("E3").Value = "No"

# Row 4
$ws7.Range("B4").Value = 23.33
$ws7.Range("C4").Value = 1
$ws7.Range("D4").Value = "No"
$ws7.Range("E4").Value = "Yes"
$ws7.Range("G4").Value = "Yes"

# Row 5
$ws7.Range("B5").Value = 39
$ws7.Range("C5").Value = 102
$ws7.Range("E5").Value = "yes"
$ws7.Range("F5").Value = "No"

# Row 6
$ws7.Range("B6").Value = 50
$ws7.Range("C6").Value = 199
$ws7.Range("D6").Value = "Yes"
$ws7.Range("E6").Value = "No"

# Row 7
$ws7.Range("B7").Value = 28
$ws7.Range("F7").Value = "Yes"

# New row 8 — copy row 7's formatting down first so the new cells pick
# up style index 7 (same as every other data row), then overwrite the
# values/text brought in by the copy.
$ws7.Range("A7:G7").Copy()
$ws7.Range("A8:G8").PasteSpecial(-4122)
$ws7.Range("A8").Value = "Auto 7"
$ws7.Range("B8").Value = 66.66
$ws7.Range("C8").Value = 160
$ws7.Range("D8").Value = "Yes"
$ws7.Range("E8").Value = "yes"
$ws7.Range("F8").Value = "Yes"
$ws7.Range("G8").Value = "No"

# Discount sheet keeps its own selection at B4, but this sheet is no
# longer the one left active when the file is saved (see below).
$ws7.Range("B4").Select()

# =====================================================================
# Bill Configuration sheet (sheet1.xml) — toggle row 17's ON/OFF cell.
# =====================================================================
$ws1 = $wb.Worksheets.Item("Bill Configuration")
$ws1.Range("B17").Value = "ON"

# Bill Configuration becomes the active sheet/selection at save time.
$ws1.Activate()
$ws1.Range("I17").Select()

$wb.Save()
